# Generate Report for Handoff
# Applies the localization-status report refresh:
#  - Status flips from "Handed back: in sync with en-US" -> "Ready for handoff"
#  - Translation method flips from "ht" -> "mt"
#  - Handoff/handback timestamps are refreshed
#  - a2160a50-*.md row now reports a stale-handback-version error
#  - A couple of report columns are narrowed / widened to fit the new content

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"
$oldMethod = "ht"
$newMethod = "mt"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e5fafeb7896e9a20ae446ba1c6c9cdec4aa0175c/e2e/a2160a50-b543-48d5-b194-3f1d31dfe14b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9123f58888fb761b1c2db98632f15ea7e867607b/e2e/a2160a50-b543-48d5-b194-3f1d31dfe14b.md."

# ColumnWidth values chosen so the resulting stored (OOXML) column width lands
# on the engine's nearest representable value to the target widths:
#   target 17.2159881591797 -> stored 17.166666666666668 (ColumnWidth 16.333333333333332)
#   target 40                -> stored 40                (ColumnWidth 39.166666666666664)
$narrowColWidth = 16.333333333333332
$wideColWidth   = 39.166666666666664

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus
$ovw.Range("E3").Value = $newStatus
$ovw.Range("F3").Value = $newStatus

$ovw.Range("G2").Value = "2016-09-07 15:13:46"
$ovw.Range("G3").Value = "2016-09-07 15:13:46"

$ovw.Columns.Item(5).ColumnWidth = $narrowColWidth
$ovw.Columns.Item(6).ColumnWidth = $narrowColWidth

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("E2").Value = $newMethod
$zhcn.Range("E3").Value = $newMethod

$zhcn.Range("H2").Value = "2016-09-07 15:13:31"
$zhcn.Range("H3").Value = "2016-09-07 15:13:31"

$zhcn.Range("P3").Value = $errorDetail

$zhcn.Columns.Item(3).ColumnWidth = $narrowColWidth
$zhcn.Columns.Item(16).ColumnWidth = $wideColWidth

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("E2").Value = $newMethod
$dede.Range("E3").Value = $newMethod

$dede.Range("H2").Value = "2016-09-07 15:13:46"
$dede.Range("H3").Value = "2016-09-07 15:13:46"

$dede.Range("P3").Value = $errorDetail

$dede.Columns.Item(3).ColumnWidth = $narrowColWidth
$dede.Columns.Item(16).ColumnWidth = $wideColWidth
